$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 143
$ws.Range("F4").Value = 132
$ws.Range("F5").Value = 1939
$ws.Range("F7").Value = 4040
$ws.Range("F8").Value = 525
$ws.Range("F9").Value = 1039
$ws.Range("F10").Value = 1306
$ws.Range("F11").Value = 653
$ws.Range("F12").Value = 365
$ws.Range("F14").Value = 2170
$ws.Range("F15").Value = 390
$ws.Range("F16").Value = 649479
$ws.Range("F17").Value = 1602
$ws.Range("F18").Value = 481
$ws.Range("F19").Value = 1429
$ws.Range("F21").Value = 538
$ws.Range("F22").Value = 1252
$ws.Range("F23").Value = 2165
$ws.Range("F24").Value = 1110
$ws.Range("F25").Value = 2667
$ws.Range("F26").Value = 1534
$ws.Range("F27").Value = 762
$ws.Range("F28").Value = 1505
$ws.Range("F29").Value = 519
$ws.Range("F30").Value = 1073
$ws.Range("F31").Value = 1074
$ws.Range("F32").Value = 74
$ws.Range("F33").Value = 2000
$ws.Range("F34").Value = 1334
$ws.Range("F35").Value = 1206
$ws.Range("F36").Value = 2051
$ws.Range("F37").Value = 1129
$ws.Range("F38").Value = 38
$ws.Range("F40").Value = 44
$ws.Range("F41").Value = 2546
$ws.Range("F42").Value = 200
$ws.Range("F43").Value = 971
$ws.Range("F44").Value = 24
$ws.Range("F45").Value = 874
$ws.Range("F46").Value = 137
$ws.Range("F47").Value = 668

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 94
$ws.Range("F10").Value = 469
$ws.Range("F11").Value = 144374
$ws.Range("F12").Value = 144374
$ws.Range("F18").Value = 223
$ws.Range("F19").Value = 330
$ws.Range("F21").Value = 401
$ws.Range("F22").Value = 401
$ws.Range("F23").Value = 113
$ws.Range("F24").Value = 77
$ws.Range("F27").Value = 523
$ws.Range("F32").Value = 319
$ws.Range("F33").Value = 266
$ws.Range("F40").Value = 183

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3112
$ws.Range("F6").Value = 234
$ws.Range("F8").Value = 816
$ws.Range("F9").Value = 1148
$ws.Range("F10").Value = 630
$ws.Range("F11").Value = 1576
$ws.Range("F12").Value = 471
$ws.Range("F13").Value = 55
$ws.Range("F14").Value = 1818

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 816
$ws.Range("F3").Value = 630
$ws.Range("F4").Value = 143
$ws.Range("F5").Value = 1576
$ws.Range("F6").Value = 471
$ws.Range("F7").Value = 132
$ws.Range("F8").Value = 1818
$ws.Range("F9").Value = 4040
$ws.Range("F11").Value = 525
$ws.Range("F12").Value = 1306
$ws.Range("F13").Value = 653
$ws.Range("F14").Value = 365
$ws.Range("F15").Value = 2170
$ws.Range("F17").Value = 390
$ws.Range("F18").Value = 649488
$ws.Range("F19").Value = 94
$ws.Range("F20").Value = 469
$ws.Range("F21").Value = 1602
$ws.Range("F22").Value = 144374
$ws.Range("F23").Value = 482
$ws.Range("F24").Value = 1429
$ws.Range("F26").Value = 538
$ws.Range("F27").Value = 1252
$ws.Range("F28").Value = 2165
$ws.Range("F29").Value = 1110
$ws.Range("F30").Value = 2667
$ws.Range("F31").Value = 1534
$ws.Range("F33").Value = 1505
$ws.Range("F34").Value = 401
$ws.Range("F35").Value = 519
$ws.Range("F36").Value = 113
$ws.Range("F37").Value = 1073
$ws.Range("F38").Value = 1074
$ws.Range("F39").Value = 77
$ws.Range("F40").Value = 74
$ws.Range("F41").Value = 2000
$ws.Range("F42").Value = 1334
$ws.Range("F43").Value = 1206
$ws.Range("F44").Value = 2051
$ws.Range("F45").Value = 1129
$ws.Range("F46").Value = 319
$ws.Range("F47").Value = 319
$ws.Range("F48").Value = 2546
$ws.Range("F49").Value = 200
$ws.Range("F50").Value = 971
$ws.Range("F51").Value = 137
$ws.Range("F52").Value = 668
